# Update the "efficient demonstration steps" results table on slide 6.
# Fills in the previously-blank KICMistral 8X7B-50 columns (MRR, Hits@1,
# Hits@3, Hits@10) for eff_demon_step rows 4, 5 and 7, and moves the
# "best MRR" bold highlight from step 10 (0.3523) to step 7 (0.3748).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$tbl = $s.Shapes.Item(8).Table

# --- eff_demon_step = 4 (table row 6) ---
$tbl.Cell(6, 2).Shape.TextFrame.TextRange.Text = "0.3229"

$tr = $tbl.Cell(6, 3).Shape.TextFrame.TextRange
$tr.Text = "0.2100"
$tr.Font.Bold = 0

$tbl.Cell(6, 4).Shape.TextFrame.TextRange.Text = "0.4200"
$tbl.Cell(6, 5).Shape.TextFrame.TextRange.Text = "0.5600"

# --- eff_demon_step = 5 (table row 7) ---
$tbl.Cell(7, 2).Shape.TextFrame.TextRange.Text = "0.3336"

$tr = $tbl.Cell(7, 3).Shape.TextFrame.TextRange
$tr.Text = "0.2300"
$tr.Font.Bold = 0

$tbl.Cell(7, 4).Shape.TextFrame.TextRange.Text = "0.3900"
$tbl.Cell(7, 5).Shape.TextFrame.TextRange.Text = "0.5400"

# --- eff_demon_step = 7 (table row 9) ---
$tr = $tbl.Cell(9, 2).Shape.TextFrame.TextRange
$tr.Text = "0.3748"
$tr.Font.Bold = 1

$tr = $tbl.Cell(9, 3).Shape.TextFrame.TextRange
$tr.Text = "0.2600"
$tr.Font.Bold = 0

$tbl.Cell(9, 4).Shape.TextFrame.TextRange.Text = "0.4600"
$tbl.Cell(9, 5).Shape.TextFrame.TextRange.Text = "0.6200"

# --- eff_demon_step = 10 (table row 12): no longer the best MRR ---
$tbl.Cell(12, 2).Shape.TextFrame.TextRange.Font.Bold = 0
